$d = $word.ActiveDocument

# Locate every "---" immediately followed by two manual line breaks
# (the "Utilities Req'd: ---" line, followed by a blank line, followed by
# the "See plans..." sentence). The document models each visual line as a
# <w:r><w:br/></w:r> run inside one giant paragraph; we need to turn the
# first of the two breaks into a real paragraph boundary while leaving the
# second break (and everything else) untouched.
$breakChar = [char]11
$pattern = "---" + $breakChar + $breakChar

$positions = New-Object System.Collections.ArrayList
$full = $d.Content.Text
$searchFrom = 0
while ($true) {
    $i = $full.IndexOf($pattern, $searchFrom)
    if ($i -lt 0) { break }
    $positions.Add($i) | Out-Null
    $searchFrom = $i + 1
}

# Walk matches back-to-front so earlier offsets stay valid as we edit.
for ($k = $positions.Count - 1; $k -ge 0; $k--) {
    $matchStart = $positions[$k]
    $firstBreak = $matchStart + 3
    $secondBreak = $matchStart + 4

    # Insert a new paragraph mark right before the second break run; this
    # keeps the second <w:r><w:br/></w:r> run intact as the first run of
    # the newly created paragraph.
    $secondBreakRange = $d.Range($secondBreak, $secondBreak + 1)
    $secondBreakRange.InsertParagraphBefore()

    # The first break run now sits alone at the end of the original
    # paragraph (just before the new paragraph mark) -- remove it so the
    # split lands exactly where the first break used to be, per the diff.
    $leftoverBreak = $d.Range($firstBreak, $firstBreak + 1)
    $leftoverBreak.Delete()
}
